{"js": "// Replace the multiplication problems in the practice-sheet table with a\n// new set of problems, as described by the commit diff.\nconst replacements = [\n  [\"83\u00d711=\", \"86\u00d749=\"],\n  [\"69\u00d714=\", \"17\u00d777=\"],\n  [\"40\u00d751=\", \"16\u00d748=\"],\n  [\"25\u00d798=\", \"60\u00d732=\"],\n  [\"28\u00d780=\", \"46\u00d779=\"],\n  [\"41\u00d740=\", \"81\u00d773=\"],\n  [\"17\u00d713=\", \"63\u00d788=\"],\n  [\"65\u00d758=\", \"76\u00d795=\"],\n  [\"88\u00d726=\", \"57\u00d746=\"],\n  [\"92\u00d719=\", \"48\u00d723=\"],\n  [\"67\u00d725=\", \"20\u00d769=\"],\n  [\"38\u00d731=\", \"21\u00d716=\"],\n  [\"98\u00d784=\", \"53\u00d713=\"],\n  [\"93\u00d797=\", \"57\u00d770=\"],\n  [\"38\u00d729=\", \"67\u00d753=\"],\n  [\"88\u00d742=\", \"63\u00d736=\"],\n  [\"75\u00d755=\", \"74\u00d784=\"],\n  [\"90\u00d719=\", \"98\u00d757=\"],\n  [\"64\u00d769=\", \"72\u00d736=\"],\n  [\"72\u00d718=\", \"11\u00d796=\"],\n  [\"88\u00d747=\", \"49\u00d727=\"],\n  [\"60\u00d731=\", \"79\u00d719=\"],\n  [\"52\u00d730=\", \"14\u00d728=\"],\n  [\"58\u00d730=\", \"25\u00d779=\"],\n  [\"51\u00d784=\", \"99\u00d722=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the practice-sheet table with a\n# new set of problems, as described by the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"83\u00d711=\", \"86\u00d749=\"),\n    @(\"69\u00d714=\", \"17\u00d777=\"),\n    @(\"40\u00d751=\", \"16\u00d748=\"),\n    @(\"25\u00d798=\", \"60\u00d732=\"),\n    @(\"28\u00d780=\", \"46\u00d779=\"),\n    @(\"41\u00d740=\", \"81\u00d773=\"),\n    @(\"17\u00d713=\", \"63\u00d788=\"),\n    @(\"65\u00d758=\", \"76\u00d795=\"),\n    @(\"88\u00d726=\", \"57\u00d746=\"),\n    @(\"92\u00d719=\", \"48\u00d723=\"),\n    @(\"67\u00d725=\", \"20\u00d769=\"),\n    @(\"38\u00d731=\", \"21\u00d716=\"),\n    @(\"98\u00d784=\", \"53\u00d713=\"),\n    @(\"93\u00d797=\", \"57\u00d770=\"),\n    @(\"38\u00d729=\", \"67\u00d753=\"),\n    @(\"88\u00d742=\", \"63\u00d736=\"),\n    @(\"75\u00d755=\", \"74\u00d784=\"),\n    @(\"90\u00d719=\", \"98\u00d757=\"),\n    @(\"64\u00d769=\", \"72\u00d736=\"),\n    @(\"72\u00d718=\", \"11\u00d796=\"),\n    @(\"88\u00d747=\", \"49\u00d727=\"),\n    @(\"60\u00d731=\", \"79\u00d719=\"),\n    @(\"52\u00d730=\", \"14\u00d728=\"),\n    @(\"58\u00d730=\", \"25\u00d779=\"),\n    @(\"51\u00d784=\", \"99\u00d722=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
